$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header values P1 and Q1, copying header style/format from O1
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# Data rows 2-25: updated values for existing columns B,C,D,E,G,H,I,O
# plus brand-new columns P (zeros) and Q (new rightmost series)
$data = @(
    @(3.461049437754639, 1.012273344699565, 0.2010353516683807, 1.346589557846372, 0.0007909771412192437, 0.0103639600143904, 0.003003242120124039, 0, 0, 2.463728656002957),
    @(3.012219475459233, 0.8887778461804601, 0.1742357590228067, 1.169434685896974, 0.0007954484576243473, 0.006955960908331782, 0.001413122548911794, 0, 0, 2.270302517013079),
    @(2.736764943485241, 0.8136252633837273, 0.1579503725555185, 1.061207587506829, 0.0007982771020863183, 0.005147121927015941, 0.0008079613893019477, 0, 0, 2.153977629180417),
    @(2.624480411616844, 0.7847226178253948, 0.1513500587162326, 1.017192682888179, 0.0007994577402732178, 0.004473790769335739, 0.0007092465621028055, 0, 0, 2.102535670372163),
    @(2.605795165316977, 0.7818589551000343, 0.1502546976355319, 1.009864846888377, 0.0007996630428404039, 0.004363709485393108, 0.0007807793246348282, 0, 0, 2.088429120335775),
    @(2.735139804791913, 0.8185120424572858, 0.1578569517526063, 1.060544225469044, 0.0007983142593263321, 0.005131682061788601, 0.001015920231903777, 0, 0, 2.137975931000625),
    @(3.306082549221571, 0.9766016279555743, 0.1917497521032914, 1.28527834738118, 0.0007925283600474695, 0.009116239639820112, 0.002614871040475641, 0, 0, 2.376163879564899),
    @(4.428704043313871, 1.28326224812389, 0.2598280663820987, 1.731834465174558, 0.0007817844010500109, 0.01947490422076714, 0.008746490764289128, 0, 0, 2.893649592927943),
    @(5.256449941141739, 1.515067829262591, 0.3111150821903124, 1.962569824270545, 0.0007744755883196566, 0.02824474436014679, 0.01566676196465444, 0, 0, 3.209534414026109),
    @(5.622536226086822, 1.62007111989135, 0.3339967342870125, 1.273846939128248, 0.0007730125263851012, 0.04373863273984924, 0.01796882537458622, 0, 0, 2.70408394381451),
    @(5.757222030316541, 1.651485538468421, 0.3424096006076525, 0.7764370143502504, 0.0007730612440655094, 0.07945867259621053, 0.01810715249733796, 0, 0, 2.269168396181328),
    @(5.7200852657719, 1.636262764454329, 0.3400106681834814, 0.3982300720573164, 0.0007742999722874511, 0.1320969609555647, 0.01683727002995195, 0, 0, 1.840233392357362),
    @(5.616492007218767, 1.606389617189109, 0.3334610125628927, 0.2062422667850683, 0.0007756840794672546, 0.1793365243734257, 0.01545973837930692, 0, 0, 1.546691373300604),
    @(5.55452199378658, 1.59138801809371, 0.3295639558311763, 0.1683537186439281, 0.0007763090466090511, 0.1911226116745723, 0.01492345635943249, 0, 0, 1.468661440270438),
    @(5.204681151289719, 1.497803298526208, 0.3076910403155182, 0.1627991411191019, 0.0007791528475942615, 0.1756743716520219, 0.01222722875592996, 0, 0, 1.421027209788861),
    @(4.992036543776351, 1.442251969277947, 0.2944834147532731, 0.2463881333921378, 0.000780624160332281, 0.1368493632530914, 0.0109422257814904, 0, 0, 1.535420185012271),
    @(4.872718468362677, 1.408334308141491, 0.2871174762456121, 0.4795118759635599, 0.0007809532167095679, 0.0846973023002846, 0.01039142100668755, 0, 0, 1.829005829348631),
    @(4.837068446173532, 1.405461151727991, 0.2849570070342367, 0.9143710867217294, 0.0007802056804552484, 0.04226461605823317, 0.01099692884226311, 0, 0, 2.25378353235584),
    @(5.038195435140381, 1.470912631290901, 0.2975143932035849, 1.89715308233373, 0.0007764355636842007, 0.0257207515521789, 0.01425705621560258, 0, 0, 3.074328429974969),
    @(5.677747994483184, 1.649263905189684, 0.3375689671483286, 2.235524352148502, 0.0007705704111492678, 0.03442693425868537, 0.0206971929507711, 0, 0, 3.451691982353623),
    @(6.09702536949959, 1.759374626196461, 0.3640934771294582, 2.406250437458837, 0.0007668821670019548, 0.04022807822875318, 0.02521435795573357, 0, 0, 3.685352499981889),
    @(5.873259693076591, 1.693939003757805, 0.3499082279727475, 2.315064427679815, 0.0007688266353109025, 0.0370987168953123, 0.02251674357046607, 0, 0, 3.578326346531583),
    @(5.029102131263528, 1.459050494817291, 0.2969496597197292, 1.973055136233995, 0.0007763251985198657, 0.02623129859828222, 0.01392874286356083, 0, 0, 3.155720742652107),
    @(4.124190613918927, 1.209311509883207, 0.2411898498322813, 1.610068420514111, 0.0007846594611968367, 0.01634323219939326, 0.007019526872079851, 0, 0, 2.722867302746863)
)

$cols = @(2, 3, 4, 5, 7, 8, 9, 15, 16, 17)  # B,C,D,E,G,H,I,O,P,Q column indices
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Cells.Item($excelRow, $cols[$c]).Value = $rowVals[$c]
    }
}
